$d = $word.ActiveDocument

# Helper: perform a Find & Replace across the whole document content.
function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: find failed for [$find]"
    }
    return $ok
}

# 1) "Assets and " + "Returns" + " " -> merge into a single run "Assets and Returns "
Replace-Text "Assets and Returns " "Assets and Returns "

# 2) "Child Education (Drawing Class)" -> "Child Education (Exam & Drawing Class)"
Replace-Text "Child Education (Drawing Class)" "Child Education (Exam & Drawing Class)"

# 3) "Child Education (Home Task)" -> "Child Education (Teacher)"
Replace-Text "Child Education (Home Task)" "Child Education (Teacher)"

# 4) "Child Education (Review)" -> "Child Education (Off)"
Replace-Text "Child Education (Review)" "Child Education (Off)"

# 5) "TIME: 4Hs" -> "TIME: 3Hs"
Replace-Text "TIME: 4Hs" "TIME: 3Hs"

# 6) "Child Education (Teacher)" -> "Child Education (Home Task)"
Replace-Text "Child Education (Teacher)" "Child Education (Home Task)"

# 7) "Child Education (Exam)" -> "Child Education (Off)"
Replace-Text "Child Education (Exam)" "Child Education (Off)"

# 8) "Law of Bangladesh" -> "Governance and International Relation"
Replace-Text "Law of Bangladesh" "Governance and International Relation"

# 9) "Governance and International Relation " -> "Law of Bangladesh and Human Rights "
Replace-Text "Governance and International Relation " "Law of Bangladesh and Human Rights "

# 10) "History and Psychology   " (3 trailing spaces across two runs) -> merge into a single run
Replace-Text "History and Psychology   " "History and Psychology   "

# 11) "Practice and Manipulate" (two runs incl. proofErr wrap) -> merge into a single run
Replace-Text "Practice and Manipulate" "Practice and Manipulate"
